# The "6 cylinder" group of data (previously on row 6, spanning the
# merged A6:A7 cell) is moved up to row 5, and the "4 cylinder" group
# (previously on row 5) moves down to row 6. Row 7 keeps its own
# statistics but the merge over A6:A7 is removed, so A7 now carries its
# own "6" label instead of being blank/merged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 5 values: old "6 cylinder" row ---
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 110
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 2.7475
$ws.Range("G5").Value = 0.1803122292025695

# --- New row 6 values: old "4 cylinder" row ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 83.33333333333333
$ws.Range("E6").Value = 18.50225211517056
$ws.Range("F6").Value = 2.886666666666667
$ws.Range("G6").Value = 0.4911551010967242

# The A6:A7 merge (with its vertical-top alignment) is no longer needed
# once A7 carries its own label, so split it apart first -- writing to
# A7 while it is still merged with A6 would redirect the write to the
# merged range's anchor cell (A6) instead.
$ws.Range("A6:A7").UnMerge()

# --- Row 7 keeps its own B:G values; only A7 changes (gets a value) ---
$ws.Range("A7").Value = 6

# A6 and A7 should look like the rest of column A (same formatting as
# A5 / A8), not the old merged-cell's vertical-top alignment style, so
# copy A5's format (not its value) onto both.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
